$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing header C1
$ws.Range("C1").Value = "Link"

# Add new headers D1:I1 (copy style from C1 so bold/border/center formatting carries over)
$ws.Range("D1").Value = "Phone Number"
$ws.Range("E1").Value = "Email"
$ws.Range("F1").Value = "Facebook Link"
$ws.Range("G1").Value = "Instagram Link"
$ws.Range("H1").Value = "Twitter Link"
$ws.Range("I1").Value = "YouTube Link"

$ws.Range("C1").Copy()
$ws.Range("D1:I1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 2 data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "+355684095232"
$ws.Range("E2").Value = "info@agroblend.com"
$ws.Range("F2").Value = "https://www.facebook.com/diamondpasta"
$ws.Range("G2").Value = "https://www.instagram.com/pastadiamond/"
$ws.Range("H2").Value = "Not found"
$ws.Range("I2").Value = "https://www.youtube.com/channel/UCyuGDc-zc4j4NmXqLvXQR7g"

# Row 3 data
$ws.Range("D3").Value = "+355 692070014"
$ws.Range("E3").Value = "info@kraco.al"
$ws.Range("F3").Value = "https://www.facebook.com/KracoAL"
$ws.Range("G3").Value = "https://www.instagram.com/kraco_nature/"
$ws.Range("H3").Value = "Not found"
$ws.Range("I3").Value = "https://www.youtube.com/channel/UCLVWuZp1esG6f4ui66vKWdg"
